$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 100
$ws.Range("B5").Value = 1

$ws.Range("C15").Select() | Out-Null
